$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geolocation")

# Rename header "grid_node" -> "region" (B3)
$ws.Cells.Item(3, 2).Value = "region"

# New grid-point rows (160-184) below the existing "geolocation" table,
# following the same alternating B27/28 - B29/30 banded-row style as the
# rows directly above them. Copy the format of the last two existing
# rows (158-159) down across the 25 new rows first, then fill in values.
$fmtSrc = $ws.Range("B158:D159")
$fmtDst = $ws.Range("B160:D184")
$fmtSrc.Copy($fmtDst)

$names = @("rez_CHE_0","rez_CHE_1","rez_CHE_10","rez_CHE_11","rez_CHE_12","rez_CHE_13","rez_CHE_14","rez_CHE_15","rez_CHE_17","rez_CHE_18","rez_CHE_19","rez_CHE_2","rez_CHE_20","rez_CHE_21","rez_CHE_22","rez_CHE_23","rez_CHE_24","rez_CHE_25","rez_CHE_3","rez_CHE_4","rez_CHE_5","rez_CHE_6","rez_CHE_7","rez_CHE_8","rez_CHE_9")
$lats  = @(46.138009727206594,46.587669909068872,47.486990272793413,46.138009727206594,46.587669909068872,47.037330090931135,47.486990272793413,46.138009727206594,47.037330090931135,47.486990272793413,46.138009727206594,47.037330090931135,46.587669909068872,47.037330090931135,47.486990272793413,46.138009727206594,46.587669909068872,47.037330090931135,46.138009727206594,46.587669909068872,47.037330090931135,47.486990272793413,46.138009727206594,46.587669909068872,47.037330090931135)
$lngs  = @(6.289697641356212,6.289697641356212,7.5802325471187393,8.225500000000002,8.225500000000002,8.225500000000002,8.225500000000002,8.8707674528812674,8.8707674528812674,8.8707674528812674,9.5160349057625293,6.289697641356212,9.5160349057625293,9.5160349057625293,9.5160349057625293,10.161302358643791,10.161302358643791,10.161302358643791,6.9349650942374748,6.9349650942374748,6.9349650942374748,6.9349650942374748,7.5802325471187393,7.5802325471187393,7.5802325471187393)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 160 + $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $lats[$i]
    $ws.Cells.Item($r, 4).Value = $lngs[$i]
}
